# Generate Report for Archive
# - Update the "Ready for handoff" status text to "In Translation" everywhere it appears.
# - Shrink the width of the now-narrower "Status" columns to match the new text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
# Target raw column width (OOXML <col width=.../>) is 13.4101845877511 characters.
# The ColumnWidth COM setter here quantizes to whole-pixel (1/6 character) steps,
# so 12.5 is the input that lands the stored width on the closest reachable step
# (13.333333333333334, i.e. pixel bucket 80/6).
$newWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
            if ("$($cell.Value2)" -eq $oldStatus) {
                $cell.Value = $newStatus
                $ws.Columns.Item($colOffset + $c).ColumnWidth = $newWidth
            }
        }
    }
}
